$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "order" column header
$ws.Range("I1").Value = "order"

# Update image_file_name (column H) for several targets that previously
# used the generic "default_target.png" placeholder, and set the new
# "order" column (column I) values for every data row.

# Row 2: cloud
$ws.Range("I2").Value = 0

# Row 3: k8s
$ws.Range("I3").Value = 4

# Row 4: k8s_jetsonnano
$ws.Range("I4").Value = 0

# Row 5: pcweb
$ws.Range("I5").Value = 0

# Row 6: pc
$ws.Range("I6").Value = 5

# Row 7: Jetson-orin
$ws.Range("I7").Value = 7

# Row 8: Jetson-xavier
$ws.Range("I8").Value = 8

# Row 9: Jetson-nano
$ws.Range("I9").Value = 9

# Row 10: s22
$ws.Range("I10").Value = 11

# Row 11: s23
$ws.Range("H11").Value = "s23.jpg"
$ws.Range("I11").Value = 10

# Row 12: Odroid-n2
$ws.Range("I12").Value = 13

# Row 13: Odroid M1
$ws.Range("H13").Value = "OdroidM1.jpg"
$ws.Range("I13").Value = 14

# Row 14: Rasberry Pi5
$ws.Range("H14").Value = "RasberryPi5.jpg"
$ws.Range("I14").Value = 12

# Row 15: Comma 3X
$ws.Range("H15").Value = "comma-3x.jpg"
$ws.Range("I15").Value = 6

# Row 16: KT cloud
$ws.Range("H16").Value = "Kt_cloud.png"
$ws.Range("I16").Value = 3

# Row 17: Amazon Web Services
$ws.Range("H17").Value = "aws.png"
$ws.Range("I17").Value = 1

# Row 18: Google Cloud Platform
$ws.Range("H18").Value = "GCP.png"
$ws.Range("I18").Value = 2

# Column widths: widen H, add width for the new I column
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 21.18

# Match the active selection recorded after the edit
[void]$ws.Range("I26").Select()
